# #5: property boat&car done
# Rebuild the 汽車 (car) sheet header row and append the missing
# property/legislator metadata columns (H:N) to the existing car rows,
# mirroring the structure already used on the 土地 / 建物 sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1: proper column header labels (B1:N1) ---
$ws.Cells.Item(1, 2).Value  = "name"
$ws.Cells.Item(1, 3).Value  = "capacity"
$ws.Cells.Item(1, 4).Value  = "owner"
$ws.Cells.Item(1, 5).Value  = "register_date"
$ws.Cells.Item(1, 6).Value  = "register_reason"
$ws.Cells.Item(1, 7).Value  = "acquire_value"
$ws.Cells.Item(1, 8).Value  = "property_category"
$ws.Cells.Item(1, 9).Value  = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- Row 2 (index 33, 轎車) : add property_category .. index columns ---
$ws.Cells.Item(2, 8).Value  = "land"
$ws.Cells.Item(2, 9).Value  = "normal"
$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = "2011-11-24"
$ws.Cells.Item(2, 11).Value = "廖國棟"
$ws.Cells.Item(2, 12).Value = 962
$ws.Cells.Item(2, 13).Value = "tmp55951"
$ws.Cells.Item(2, 14).Value = 33

# --- Row 3 (index 34, 休旅車) : add property_category .. index columns ---
$ws.Cells.Item(3, 8).Value  = "land"
$ws.Cells.Item(3, 9).Value  = "normal"
$ws.Cells.Item(3, 10).NumberFormat = "@"
$ws.Cells.Item(3, 10).Value = "2011-11-24"
$ws.Cells.Item(3, 11).Value = "廖國棟"
$ws.Cells.Item(3, 12).Value = 962
$ws.Cells.Item(3, 13).Value = "tmp55951"
$ws.Cells.Item(3, 14).Value = 34

# --- carry over the bold / bordered / centered header formatting from the
#     existing header cell (B1) onto the newly added header cells (H1:N1) ---
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false


